$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Femacal de La Calera - Arveja Verde".
# It becomes the new row 66, pushing the existing rows 66-89 down to 67-90.
$ws.Rows.Item(66).Insert()

# Fill in the newly inserted row 66 with the new record's values.
$ws.Range("A66").Value = 3
$ws.Range("B66").Value = "Femacal de La Calera"
$ws.Range("C66").Value = "Coquimbo"
$ws.Range("D66").Value = 45119
$ws.Range("E66").Value = 5
$ws.Range("F66").Value = 100112022
$ws.Range("G66").Value = "Arveja Verde"
$ws.Range("H66").Value = "Perfection"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 30
$ws.Range("K66").Value = 24000
$ws.Range("L66").Value = 24000
$ws.Range("M66").Value = 24000
$ws.Range("N66").Value = "$/saco 25 kilos"
$ws.Range("O66").Value = "Provincia de Limarí"
$ws.Range("P66").Value = 960
$ws.Range("Q66").Value = 25
$ws.Range("R66").Value = "Hortaliza"
